$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------
# Part 1: the existing "Meeting" row (18/11/2021, row 11) has a w:shd
# on its 3rd cell that still carries theme-fill metadata
# (w:themeFill="accent4" w:themeFillTint="99"); the target keeps the
# same color but as a plain w:fill="FFD966" with no theme reference.
#
# Setting BackgroundPatternColor directly on a cell that already has
# theme attributes always re-derives/keeps the theme reference, so it
# cannot be stripped that way. Instead, insert a fresh row cloned from
# the row immediately above (which has no theme fill on that column),
# move the old row's text into the clone, fix the one color that
# differs, then delete the original themed row.
# ---------------------------------------------------------------------

$oldRow = $t.Rows.Item(11)

$oldText1 = $oldRow.Cells.Item(1).Range.Text
$oldText2 = $oldRow.Cells.Item(2).Range.Text
$oldText3 = $oldRow.Cells.Item(3).Range.Text
$oldText4 = $oldRow.Cells.Item(4).Range.Text

# Inserting before row 11 clones the (theme-free) formatting of row 10.
$null = $t.Rows.Add($oldRow)

$newRow = $t.Rows.Item(11)
$srcRow = $t.Rows.Item(12)

$newRow.Cells.Item(1).Range.Text = $oldText1
$newRow.Cells.Item(2).Range.Text = $oldText2
$newRow.Cells.Item(3).Range.Text = $oldText3
$newRow.Cells.Item(4).Range.Text = $oldText4

# Column 3 ("type") is the only one whose fill differs between row 10
# (green, 7AB648) and row 11 (gold, FFD966) - push the correct color.
$newRow.Cells.Item(3).Shading.BackgroundPatternColor = 6740479

# Remove the old, themed row (now pushed down to index 12).
$srcRow.Delete()

# ---------------------------------------------------------------------
# Part 2: append a brand-new row documenting the progress report.
# Rows.Add() with no argument appends after the last row, cloning the
# last row's (theme-free) formatting - exactly what the new row needs
# for 3 of its 4 cells; only the "type" cell's color must change.
# ---------------------------------------------------------------------

$addedRow = $t.Rows.Add()

$addedRow.Cells.Item(1).Range.Text = "6/12/2021"
$addedRow.Cells.Item(2).Range.Text = "2 Hours 10 minutes"
$addedRow.Cells.Item(3).Range.Text = "Progress Report"
$addedRow.Cells.Item(4).Range.Text = "Began the writeup of the progress report, demonstrating what has been achieved and how the project is on track to meet its goals."

# Column 3 ("type") differs between the last existing row (green,
# 7AB648) and this new row (gold, FFD966).
$addedRow.Cells.Item(3).Shading.BackgroundPatternColor = 6740479
